# Auto-update draw results: append the 2025-10-05 Pick 4 draw as a new
# row (19) at the bottom of the "Results" sheet, extending the used
# range from A1:E18 to A1:E19 (dimension + ignoredErrors sqref follow).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19

# Columns A (plain date text) and C (all-digit phase code) look like a
# date/number to Excel's literal-entry parser and would otherwise get
# silently coerced into a date serial / numeric value. Mark just those
# two cells as Text first so the values round-trip as the literal
# strings shown in the source feed (matches the sheet's existing
# numberStoredAsText ignored-error hint). B/D/E are never ambiguous, so
# leave their formatting untouched.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 3).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-10-05"
$ws.Cells.Item($row, 2).Value = "Pick 4"
$ws.Cells.Item($row, 3).Value = "251005"
$ws.Cells.Item($row, 4).Value = "4-7-5-5"
$ws.Cells.Item($row, 5).Value = "2025-10-05T21:34:52.878+04:00"
